$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new day's data as row 18, matching the formatting of the
# preceding row (date style on column A) by copying formats down.
$ws.Range("A17:B17").Copy()
$ws.Range("A18:B18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A18").Value = 45980
$ws.Range("B18").Value = 55

$ws.Range("A18:B18").Select()
